$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the existing header style (bold,
# bordered, centered) from H1 so the new columns look the same as B1:H1.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I34 with a constant 1, and J2:J34 mirroring the value already in
# column H (IP) for that row.
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hVal
}
